$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.445.76'
$ws.Range("E2").Value = '  -1.34%  '

# Row 3
$ws.Range("D3").Value = '3.761.52'
$ws.Range("E3").Value = '  +0.38%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").Value = '''593.87'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("D6").Value = '''165.81'
$ws.Range("E6").Value = '  -0.64%  '

# Row 7
$ws.Range("D7").Value = '3.766.83'
$ws.Range("E7").Value = '  +0.56%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("E9").Value = '  -0.17%  '

# Row 10
$ws.Range("E10").Value = '  -0.13%  '

# Row 11
$ws.Range("D11").Value = '''6.35'
$ws.Range("E11").Value = '  -1.79%  '

# Row 12
$ws.Range("D12").Value = '''0.448'
$ws.Range("E12").Value = '  +0.41%  '

# Row 13
$ws.Range("E13").Value = '  -1.64%  '

# Row 14
$ws.Range("D14").Value = '''36.10'
$ws.Range("E14").Value = '  +0.41%  '

# Row 15
$ws.Range("D15").Value = '4.388.06'
$ws.Range("E15").Value = '  +0.28%  '

# Row 16
$ws.Range("D16").Value = '3.758.67'
$ws.Range("E16").Value = '  +0.47%  '

# Row 17
$ws.Range("D17").Value = '''18.42'
$ws.Range("E17").Value = '  +2.60%  '

# Row 18
$ws.Range("D18").Value = '67.386.63'
$ws.Range("E18").Value = '  -1.48%  '

# Row 19
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$ws.Range("D20").Value = '''6.99'
$ws.Range("E20").Value = '  -0.04%  '

# Row 21
$ws.Range("D21").Value = '''10.01'
$ws.Range("E21").Value = '  -6.75%  '

# Row 22
$ws.Range("D22").Value = '''455.51'
$ws.Range("E22").Value = '  -2.06%  '

# Row 23
$ws.Range("D23").Value = '''0.696'
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("E24").Value = '  +7.03%  '

# Row 25
$ws.Range("D25").Value = '''83.22'
$ws.Range("E25").Value = '  -1.44%  '

# Row 26
$ws.Range("E26").Value = '  -1.80%  '

# Row 27
$ws.Range("D27").Value = '''11.88'
$ws.Range("E27").Value = '  -0.89%  '

# Row 28
$ws.Range("D28").Value = '''10.13'
$ws.Range("E28").Value = '  +1.21%  '

# Row 29
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("E30").Value = '  -0.11%  '

# Row 31
$ws.Range("D31").Value = '''7.27'
$ws.Range("E31").Value = '  -0.10%  '

# Row 32
$ws.Range("D32").Value = '''29.67'
$ws.Range("E32").Value = '  -0.41%  '

# Row 33
$ws.Range("E33").Value = '  +0.41%  '

# Row 34
$ws.Range("D34").Value = '''9.17'
$ws.Range("E34").Value = '  -0.12%  '

# Row 35
$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.25%  '

# Row 36
$ws.Range("D36").Value = '3.712.39'
$ws.Range("E36").Value = '  +0.28%  '

# Row 37
$ws.Range("E37").Value = '  -0.29%  '

# Row 38
$ws.Range("E38").Value = '  -1.49%  '

# Row 39
$ws.Range("E39").Value = '  -0.89%  '

# Row 40
$ws.Range("D40").Value = '''0.993'
$ws.Range("E40").Value = '  -0.50%  '

# Row 41
$ws.Range("D41").Value = '''5.74'
$ws.Range("E41").Value = '  -1.04%  '

# Row 42
$ws.Range("E42").Value = '  -0.13%  '

# Row 43
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("D44").Value = '''45.27'
$ws.Range("E44").Value = '  +2.83%  '

# Row 45
$ws.Range("E45").Value = '  -1.64%  '

# Row 46
$ws.Range("D46").Value = '''47.04'
$ws.Range("E46").Value = '  +2.47%  '

# Row 47
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = '''8.33'
$ws.Range("E47").Value = '  -2.52%  '

# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''148.61'
$ws.Range("E48").Value = '  +1.81%  '

# Row 49
$ws.Range("E49").Value = '  -4.26%  '

# Row 50
$ws.Range("D50").Value = '''389.46'
$ws.Range("E50").Value = '  +0.11%  '

# Row 51
$ws.Range("D51").Value = '''26.02'
$ws.Range("E51").Value = '  +1.69%  '
